$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("README")

# Write the new "Notes" bullet (#4.) into B13, replacing the previously blank cell.
$ws.Range("B13").Value = "#4. Enter lists using the pipe (“|”) character as a separator, for applicable DwC terms only."
$ws.Range("B13").WrapText = $true
$ws.Range("B13").Font.Size = 12
$ws.Range("B13").Font.Color = 0
$ws.Range("B13").Font.Bold = $false

# "#4." is bold, the remainder of the sentence is regular weight.
$bold = $ws.Range("B13").Characters(1, 3)
$bold.Font.Name = "Calibri"
$bold.Font.Size = 12
$bold.Font.Color = 0
$bold.Font.Bold = $true

$rest = $ws.Range("B13").Characters(4, 90)
$rest.Font.Name = "Calibri"
$rest.Font.Size = 12
$rest.Font.Color = 0
$rest.Font.Bold = $false

# Make README the active tab (previously Occurrences was active/selected).
$ws.Activate()
$ws.Range("A1").Select()
